# Append two new employee rows (108 rows total now span A1:H9) to the
# Cloud_EMS sheet, matching the existing David / Male / 9578821821 / 30 /
# 40310 (DOJ) / 50000 pattern used by the prior rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("A8").Value = 107
$ws.Range("B8").Value = "David"
$ws.Range("C8").Value = 9578821821
$ws.Range("D8").Value = 30
$ws.Range("F8").Value = "Male"
$ws.Range("G2").Copy($ws.Range("G8"))
$ws.Range("H8").Value = 50000

# --- Row 9 ---
$ws.Range("A9").Value = 108
$ws.Range("B9").Value = "David"
$ws.Range("C9").Value = 9578821821
$ws.Range("D9").Value = 30
$ws.Range("F9").Value = "Male"
$ws.Range("G2").Copy($ws.Range("G9"))
$ws.Range("H9").Value = 50000
